# "add user list to project"
#
# Adds a new "users" column to the "project hours" sheet: a header in E1
# and, for each existing data row, the (stringified) list of users
# associated with that project in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project hours")

# New header cell.
$ws.Range("E1").Value = "users"

# New per-row data.
$ws.Range("E2").Value = "['Arun Lakshmanan']"
$ws.Range("E3").Value = "['Kyle Pieper']"

# Give the new header cell (E1) the same look as the rest of the header
# row (bold font, thin border, centered/top-aligned) by copying the
# formatting from the neighboring "project" header cell, rather than
# setting Font/Borders/Alignment properties individually (which would
# create redundant, unused style entries in the workbook's style table).
$ws.Range("B1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
